$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  D=18795.86;          F=20.6387225309397},
    @{Row=3;  D=12075.86;          F=314.703159604627},
    @{Row=4;  D=12075.86;          F=315.785820599274},
    @{Row=5;  D=12075.86;          F=298.391327678514},
    @{Row=6;  D=12075.86;          F=149.290850188205},
    @{Row=7;  D=9791.86;           F=250.18380862747},
    @{Row=8;  D=9791.86;           F=377.865700216756},
    @{Row=9;  D=9791.86;           F=369.901197534505},
    @{Row=10; D=9791.86;           F=367.627171573642},
    @{Row=11; D=9791.86;           F=367.340639677825},
    @{Row=12; D=9791.86;           F=341.902574424583},
    @{Row=13; D=9791.86;           F=208.002437228832},
    @{Row=14; D=9791.86;           F=207.815067956182},
    @{Row=15; D=9791.86;           F=340.607855821684}
)

foreach ($entry in $data) {
    $ws.Cells.Item($entry.Row, 4).Value = $entry.D
    $ws.Cells.Item($entry.Row, 6).Value = $entry.F
}
